$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new column before column N (14th column) - shifts N..P to O..Q
$ws.Columns("N").Insert()

# The newly inserted column should carry the same width as column M (10 characters), not bestFit
$ws.Columns("N").ColumnWidth = 9.17

# Set selection to R8 to match the final state
$ws.Range("R8").Select()
